# Applies the updated cryptocurrency price/volume snapshot (and the two
# swapped-rank row pairs) to the coin list, matching the Nov 6 2023
# GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.499.98"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.915.13"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.724"
$ws.Range("E5").Value = "  +11.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "248.60"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.79"
$ws.Range("E8").Value = "  -3.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.357"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.79"
$ws.Range("E10").Value = "  +7.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0740"
$ws.Range("E11").Value = "  +2.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0990"
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.195.08"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "12.67"
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.919.40"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.93"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.529.41"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.10"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0833"
$ws.Range("E20").Value = "  +0.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.18"
$ws.Range("E21").Value = "  +3.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "242.39"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("E23").Value = "  +3.98%  "
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("E26").Value = "  +6.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.82"
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.69"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.135"
$ws.Range("E29").Value = "  +5.60%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.83"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.140.47"
$ws.Range("E31").Value = "  +19.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.36"
$ws.Range("E32").Value = "  +4.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.98"
$ws.Range("E33").Value = "  +13.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0580"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.26"
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.919"
$ws.Range("E37").Value = "  -3.17%  "
$ws.Range("E38").Value = "  +9.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.07"
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.82"
$ws.Range("E40").Value = "  +12.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.02"
$ws.Range("E41").Value = "  +7.23%  "
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0659"
$ws.Range("E43").Value = "  +3.48%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0211"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.49"
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.349.54"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.63"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.22"
$ws.Range("E50").Value = "  -3.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.95"
$ws.Range("E51").Value = "  -5.85%  "
